$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the PCR-protocol table. The write order below matches the order the
# values were originally typed in (first occurrences in this particular
# sequence), so that the shared-strings table comes out in the same order.
$ws.Range("A1").Value = "Temperatur"
$ws.Range("B1").Value = "Zeit"
$ws.Range("C1").Value = "Zyklen"

$ws.Range("A2").Value = "95 °C"
$ws.Range("B2").Value = "60s"
$ws.Range("C2").Value = "1x"

$ws.Range("A4").Value = "60 °C"
$ws.Range("B3").Value = "10s"
$ws.Range("B4").Value = "30s"
$ws.Range("C3").Value = "45x"

$ws.Range("A5").Value = "40 °C"

$ws.Range("A3").Value = "95 °C"
$ws.Range("B5").Value = "30s"
$ws.Range("C5").Value = "1x"

# Column A was widened to fit the temperature labels.
$ws.Columns.Item(1).ColumnWidth = 10.42578125

# Leave the selection where the author last clicked.
$ws.Range("C12").Select()
